# Weekly update: insert one new price-report row for "Pepino ensalada" at
# Mercado Mayorista Lo Valledor de Santiago, pushing the existing rows
# 580:637 down to 581:638 and growing the used range to A1:R638.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 580 (shifts rows 580:637 -> 581:638).
$ws.Rows.Item(580).Insert()

# Populate the newly inserted row 580 with this week's record.
$ws.Cells.Item(580, 1).Value  = 6
$ws.Cells.Item(580, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(580, 3).Value  = "Metropolitana"
$ws.Cells.Item(580, 4).Value  = 45124
$ws.Cells.Item(580, 5).Value  = 13
$ws.Cells.Item(580, 6).Value  = 100112043
$ws.Cells.Item(580, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(580, 8).Value  = "Sin especificar"
$ws.Cells.Item(580, 9).Value  = "Primera"
$ws.Cells.Item(580, 10).Value = 430
$ws.Cells.Item(580, 11).Value = 11000
$ws.Cells.Item(580, 12).Value = 12000
$ws.Cells.Item(580, 13).Value = 11581
$ws.Cells.Item(580, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(580, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(580, 16).Value = 193
$ws.Cells.Item(580, 17).Value = 60
$ws.Cells.Item(580, 18).Value = "Hortaliza"
